$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.404.96'
$ws.Range("E2").Value = '  -1.42%  '

$ws.Range("D3").Value = '2.222.30'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +9.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.65%  '

$ws.Range("E8").Value = '  -0.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.597'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.63'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -6.14%  '

$ws.Range("E11").Value = '  -2.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.35'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.62'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -5.71%  '

$ws.Range("E14").Value = '  +14.62%  '

$ws.Range("E15").Value = '  -1.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.88'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.08%  '

$ws.Range("D17").Value = '2.555.12'
$ws.Range("E17").Value = '  -0.84%  '

$ws.Range("D18").Value = '2.227.75'
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("D19").Value = '42.377.86'
$ws.Range("E19").Value = '  -1.49%  '

$ws.Range("E20").Value = '  -2.18%  '

$ws.Range("E21").Value = '  +5.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.32'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +14.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '235.05'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.86'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()

$ws.Range("E28").Value = '  -6.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.42'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -8.96%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.31'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.13'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.26'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("E34").Value = '  -2.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.60'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.98'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.126'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.17'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("E39").Value = '  +1.31%  '

$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("E41").Value = '  -4.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.42'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("E43").Value = '  -2.11%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.32'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -8.02%  '

$ws.Range("E46").Value = '  -2.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.33'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -5.24%  '

$ws.Range("E48").Value = '  +2.58%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.15'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.38'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.64'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +4.02%  '
